# Commit: "added checkbox for insurance. Added new tables, added insurance
# coeffs to those."
#
# This adds two new coefficient tables ("Назначение-страховка" and
# "Страховка-здоровье") to the ku_tables workbook, wired up the same way as
# the existing "Назначение-*" / "Питание-*" pairs, and nudges the selection
# on a couple of existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "Назначение-страховка" - inserted right before the sheet
#    that is currently active ("Питание-здоровье"), i.e. between
#    "Назначение-возраст" and "Питание-здоровье".
# ---------------------------------------------------------------------
$wsNaznStrah = $wb.Worksheets.Add()
$wsNaznStrah.Name = "Назначение-страховка"

$wsNaznStrah.Cells.Item(1, 1).Value = "X"
$wsNaznStrah.Cells.Item(1, 3).Value = "Есть"
$wsNaznStrah.Cells.Item(1, 2).Value = "Нет"

$naznRows = @(
    @("транзитные", 1, 1),
    @("делового назначения", 1, 1),
    @("для спокойного (в том числе семейного) отдыха", 1, 1),
    @("для активного отдыха", 0.5, 1),
    @("для занятия спортом", 0.8, 1),
    @("для укрепления и восстановления здоровья", 1, 1),
    @("для экскурсионного отдыха", 0.8, 1),
    @("для шоппинга", 1, 1)
)

$r = 2
foreach ($row in $naznRows) {
    $wsNaznStrah.Cells.Item($r, 1).Value = $row[0]
    $wsNaznStrah.Cells.Item($r, 2).Value = $row[1]
    $wsNaznStrah.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$wsNaznStrah.Columns.Item(1).ColumnWidth = 37.140625

# ---------------------------------------------------------------------
# 2) New sheet "Страховка-здоровье" - appended after the last sheet
#    ("Питание-возраст").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStrahZdor = $wb.Worksheets.Add($null, $lastSheet)
$wsStrahZdor.Name = "Страховка-здоровье"

$wsStrahZdor.Cells.Item(1, 1).Value = "X"
$wsStrahZdor.Cells.Item(1, 2).Value = "плохое"
$wsStrahZdor.Cells.Item(1, 3).Value = "ниже среднего"
$wsStrahZdor.Cells.Item(1, 4).Value = "среднее"
$wsStrahZdor.Cells.Item(1, 5).Value = "выше среднего"
$wsStrahZdor.Cells.Item(1, 6).Value = "хорошее"

$wsStrahZdor.Cells.Item(2, 1).Value = "Есть"
$wsStrahZdor.Cells.Item(2, 2).Value = 1
$wsStrahZdor.Cells.Item(2, 3).Value = 1
$wsStrahZdor.Cells.Item(2, 4).Value = 1
$wsStrahZdor.Cells.Item(2, 5).Value = 1
$wsStrahZdor.Cells.Item(2, 6).Value = 1

$wsStrahZdor.Cells.Item(3, 1).Value = "Нет"
$wsStrahZdor.Cells.Item(3, 2).Value = 0.5
$wsStrahZdor.Cells.Item(3, 3).Value = 0.7
$wsStrahZdor.Cells.Item(3, 4).Value = 1
$wsStrahZdor.Cells.Item(3, 5).Value = 1
$wsStrahZdor.Cells.Item(3, 6).Value = 1

$wsStrahZdor.Columns.Item(5).ColumnWidth = 15.42578125

# ---------------------------------------------------------------------
# 3) Selection tweaks on existing sheets.
# ---------------------------------------------------------------------
$wsNaznZdor = $wb.Worksheets.Item("Назначение-здоровье")
$wsNaznZdor.Range("A1:C9").Select()

$wsPitanieZdor = $wb.Worksheets.Item("Питание-здоровье")
$wsPitanieZdor.Range("A1:F3").Select()

$wsStrahZdor.Range("C4").Select()

# ---------------------------------------------------------------------
# 4) Final active sheet / selection: "Назначение-страховка", cell B9 -
#    matches workbookView activeTab="4" in the target workbook.
# ---------------------------------------------------------------------
$wsNaznStrah.Activate()
$wsNaznStrah.Range("B9").Select()
